$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

# Delete the row for "73431991 HERNAN DARIO MIRANDA FONSECA" (row 19)
$ws.Rows.Item(19).Delete()

# Update the "2507" periods to "2508" across the data table (E16:E19 after delete)
$ws.Range("E16:E19").Value = "2508"

# Update Valor Mora total (E11) 293440 -> 236500
$ws.Range("E11").Value = 236500

# Update Cant. Trabajadores (C13) from 5 -> 4
$ws.Range("C13").Value = 4

# Update the last data row (HYLEANA) G value from 908526 -> 1423500
$ws.Range("G19").Value = 1423500
